$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 46478.22
$ws.Range("J64").Value = 3214.0715
$ws.Range("L64").Value = 3214.0715
$ws.Range("N64").Value = -3710.0715
# Row 67
$ws.Range("H67").Value = 46478.22
$ws.Range("J67").Value = 3214.0715
$ws.Range("L67").Value = 3214.0715
$ws.Range("N67").Value = -4930.0715
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 129
$ws.Range("H129").Value = 259968.69
$ws.Range("J129").Value = 309645.72
$ws.Range("L129").Value = 928937.1599999999
$ws.Range("N129").Value = -938937.1599999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2163.775
$ws.Range("I61").Value = 1422.55
$ws.Range("J61").Value = 2905
$ws.Range("K61").Value = 1422.55
$ws.Range("L61").Value = 2905
$ws.Range("M61").Value = -1210.55
$ws.Range("N61").Value = -3329
# Row 63
$ws.Range("H63").Value = 2765
$ws.Range("I63").Value = 1702.5
$ws.Range("J63").Value = 3190
$ws.Range("K63").Value = 1702.5
$ws.Range("L63").Value = 3190
$ws.Range("M63").Value = -1016.5
$ws.Range("N63").Value = -4562
# Row 66
$ws.Range("H66").Value = 2765
$ws.Range("I66").Value = 1702.5
$ws.Range("J66").Value = 3190
$ws.Range("K66").Value = 8512.5
$ws.Range("L66").Value = 15950
$ws.Range("M66").Value = -5080.5
$ws.Range("N66").Value = -22814
# Row 88
$ws.Range("H88").Value = 3044.5557
$ws.Range("J88").Value = 2943
$ws.Range("L88").Value = 2943
$ws.Range("N88").Value = -3755
# Row 91
$ws.Range("H91").Value = 3044.5557
$ws.Range("J91").Value = 2943
$ws.Range("L91").Value = 2943
$ws.Range("N91").Value = -5751
# Row 122
$ws.Range("H122").Value = 2694.4
$ws.Range("I122").Value = 2873
$ws.Range("J122").Value = 1980
$ws.Range("K122").Value = 8619
$ws.Range("L122").Value = 5940
$ws.Range("M122").Value = -6169
$ws.Range("N122").Value = -10840
# Row 132
$ws.Range("H132").Value = 3295.4243
$ws.Range("I132").Value = 3483.1667
$ws.Range("J132").Value = 1418
$ws.Range("K132").Value = 10449.5001
$ws.Range("L132").Value = 4254
$ws.Range("M132").Value = -7919.500100000001
$ws.Range("N132").Value = -9314
# Row 136
$ws.Range("H136").Value = 2163.775
$ws.Range("I136").Value = 1422.55
$ws.Range("J136").Value = 2905
$ws.Range("K136").Value = 4267.65
$ws.Range("L136").Value = 8715
$ws.Range("M136").Value = -1717.65
$ws.Range("N136").Value = -13815

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 17540.125
$ws.Range("J82").Value = 21940.5
$ws.Range("L82").Value = 21940.5
$ws.Range("N82").Value = -22706.5
# Row 85
$ws.Range("H85").Value = 17540.125
$ws.Range("J85").Value = 21940.5
$ws.Range("L85").Value = 21940.5
$ws.Range("N85").Value = -24592.5
# Row 99
$ws.Range("H99").Value = 1859.4286
$ws.Range("I99").Value = 1785
$ws.Range("J99").Value = 1958.6666
$ws.Range("K99").Value = 1785
$ws.Range("L99").Value = 1958.6666
$ws.Range("M99").Value = -287
$ws.Range("N99").Value = -4954.6666
# Row 107
$ws.Range("H107").Value = 166733730
$ws.Range("I107").Value = 333467000
$ws.Range("J107").Value = 440
$ws.Range("K107").Value = 333467000
$ws.Range("L107").Value = 440
$ws.Range("M107").Value = -333465080
$ws.Range("N107").Value = -4280

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 444
$ws.Range("I22").Value = 321.6
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 321.6
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = 28.39999999999998
$ws.Range("N22").Value = -1450
# Row 99
$ws.Range("H99").Value = 10872.467
$ws.Range("I99").Value = 4278
$ws.Range("J99").Value = 29007.25
$ws.Range("K99").Value = 4278
$ws.Range("L99").Value = 29007.25
$ws.Range("M99").Value = -2780
$ws.Range("N99").Value = -32003.25
# Row 126
$ws.Range("H126").Value = 10872.467
$ws.Range("I126").Value = 4278
$ws.Range("J126").Value = 29007.25
$ws.Range("K126").Value = 12834
$ws.Range("L126").Value = 87021.75
$ws.Range("M126").Value = -10364
$ws.Range("N126").Value = -91961.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 842.9899
$ws.Range("J131").Value = 862.53766
$ws.Range("L131").Value = 2587.61298
$ws.Range("N131").Value = -12667.61298
# Row 138
$ws.Range("H138").Value = 2157.4443
$ws.Range("I138").Value = 1680.6666
$ws.Range("J138").Value = 3111
$ws.Range("K138").Value = 5041.9998
$ws.Range("L138").Value = 9333
$ws.Range("M138").Value = 98.0002000000004
$ws.Range("N138").Value = -19613
# Row 140
$ws.Range("H140").Value = 8539.933999999999
$ws.Range("I140").Value = 12699.889
$ws.Range("J140").Value = 2300
$ws.Range("K140").Value = 38099.667
$ws.Range("L140").Value = 6900
$ws.Range("M140").Value = -32919.667
$ws.Range("N140").Value = -17260
# Row 141
$ws.Range("H141").Value = 4007.2727
$ws.Range("J141").Value = 5450
$ws.Range("L141").Value = 16350
$ws.Range("N141").Value = -26710

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 404203.2
$ws.Range("I70").Value = 669002.7
$ws.Range("J70").Value = 7004
$ws.Range("K70").Value = 669002.7
$ws.Range("L70").Value = 7004
$ws.Range("M70").Value = -668732.7
$ws.Range("N70").Value = -7544
# Row 73
$ws.Range("H73").Value = 404203.2
$ws.Range("I73").Value = 669002.7
$ws.Range("J73").Value = 7004
$ws.Range("K73").Value = 669002.7
$ws.Range("L73").Value = 7004
$ws.Range("M73").Value = -668066.7
$ws.Range("N73").Value = -8876
# Row 97
$ws.Range("H97").Value = 166670260
$ws.Range("I97").Value = 200003900
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 200003900
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -200003404
$ws.Range("N97").Value = -2992
# Row 122
$ws.Range("H122").Value = 2856
$ws.Range("I122").Value = 2998.6667
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8996.000100000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6546.000100000001
$ws.Range("N122").Value = -10900

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1742.0526
$ws.Range("J100").Value = 1812.4375
$ws.Range("L100").Value = 1812.4375
$ws.Range("N100").Value = -2894.4375
# Row 136
$ws.Range("H136").Value = 2087.4167
$ws.Range("I136").Value = 1913.5454
$ws.Range("K136").Value = 5740.6362
$ws.Range("M136").Value = -3190.6362

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 645.7059
$ws.Range("I136").Value = 520.2619
$ws.Range("J136").Value = 1231.1111
$ws.Range("K136").Value = 1560.7857
$ws.Range("L136").Value = 3693.3333
$ws.Range("M136").Value = 989.2143000000001
$ws.Range("N136").Value = -8793.3333
# Row 138
$ws.Range("H138").Value = 65429
$ws.Range("J138").Value = 65429
$ws.Range("L138").Value = 65429
$ws.Range("N138").Value = -75709
